# Updated symbol list (price/volume refresh) - GitHub Actions scheduled run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value. Values are written with a leading apostrophe so
# Excel stores numeric-looking strings (prices) and percent-looking strings
# (volumes) as literal text, matching the sheet's existing text-formatted columns.
$updates = @{
    "D2" = "278.21"
    "E2" = "1.89%"
    "D3" = "27.25"
    "E3" = "1.69%"
    "D4" = "4.874"
    "D5" = "0.06371"
    "E5" = "0.75%"
    "D6" = "6.976"
    "E6" = "0.93%"
    "D7" = "1.256"
    "E7" = "-9.14%"
    "D8" = "0.8807"
    "E8" = "-0.41%"
    "D9" = "0.1512"
    "E9" = "2.65%"
    "D10" = "0.05136"
    "E10" = "0.74%"
    "D11" = "0.07541"
    "E11" = "2.11%"
    "D12" = "0.02963"
    "E12" = "-6.84%"
    "D13" = "0.09012"
    "E13" = "-0.39%"
    "E14" = "0.33%"
    "D15" = "0.0006416"
    "E15" = "1.62%"
    "D16" = "0.005928"
    "E16" = "-1.50%"
    "D17" = "3.464"
    "E17" = "-0.22%"
    "D18" = "3.319"
    "E18" = "-0.98%"
    "E20" = "0.93%"
    "E21" = "0.20%"
    "D22" = "3.907"
    "E22" = "-0.38%"
    "D23" = "0.04413"
    "E23" = "1.75%"
    "E24" = "-0.78%"
    "E25" = "6.62%"
    "D26" = "0.0001199"
    "E26" = "-0.34%"
    "E27" = "13.78%"
    "D40" = "0.04154"
    "E40" = "2.49%"
    "D41" = "0.006848"
    "E41" = "7.60%"
    "D42" = "0.1179"
    "E42" = "1.24%"
    "D43" = "0.002018"
    "E43" = "-9.32%"
    "D44" = "0.01122"
    "E44" = "-11.13%"
    "D45" = "0.00005184"
    "E45" = "-2.98%"
    "D46" = "1.487"
    "E46" = "-36.89%"
    "E47" = "-4.77%"
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = "'" + $updates[$cellRef]
}

# Re-normalize style back to the default "Normal" cell style; assigning a
# leading apostrophe causes Excel to apply a quote-prefix style, but these
# cells originally carried no explicit style.
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Style = "Normal"
}
